$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range (A1:I11, header row included) by column A ascending
$rng = $ws.Range("A1:I11")
$key1 = $ws.Range("A2")
$rng.Sort($key1, 1, $null, $null, 1, $null, 1, 1)

# Update the active selection/cell on the sheet
$ws.Range("M5").Select()
